$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect the new "through" date
$ws.Name = "Through 2022-11-27"

# Update the label for the November row
$ws.Range("A12").Value = "November (through 11-27)"

# Update November row values (row 12)
$ws.Range("B12").Value = 29
$ws.Range("C12").Value = 70
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = 60
$ws.Range("F12").Value = 47
$ws.Range("G12").Value = 192
$ws.Range("H12").Value = 183
$ws.Range("I12").Value = 104

# Update Total row values (row 13)
$ws.Range("B13").Value = 287
$ws.Range("C13").Value = 556
$ws.Range("D13").Value = 810
$ws.Range("E13").Value = 675
$ws.Range("F13").Value = 529
$ws.Range("G13").Value = 1249
$ws.Range("H13").Value = 1624
$ws.Range("I13").Value = 1502
